# Adds a new "2022-Q4" sheet (duplicated from "2022-Q3" structure) right
# after "总计", fills it with the new quarter's fund-holding data, and
# updates the "总计" (totals) sheet with the new quarter's summary row.

$wb = $excel.ActiveWorkbook

# Helper: set a cell's value while forcing text storage (so numeric-looking
# strings like "4.09" or fund codes like "000880" are not silently coerced
# to the Number type / lose leading zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (keeps the
#    same headers/column styles) and placing it right after "总计".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsTotal) | Out-Null

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Row 2
Set-TextValue $wsQ4.Range("B2") "000880"
Set-TextValue $wsQ4.Range("C2") "富国研究精选灵活配置混合A"
Set-TextValue $wsQ4.Range("D2") "4.09"
Set-TextValue $wsQ4.Range("E2") "93.15"
Set-TextValue $wsQ4.Range("F2") "8.33"
Set-TextValue $wsQ4.Range("G2") "0.3407"
$wsQ4.Range("H2").Value = 1

# Row 3
Set-TextValue $wsQ4.Range("B3") "016313"
Set-TextValue $wsQ4.Range("C3") "富国研究精选灵活配置混合C"
Set-TextValue $wsQ4.Range("D3") "0.16"
Set-TextValue $wsQ4.Range("E3") "93.15"
Set-TextValue $wsQ4.Range("F3") "8.33"
Set-TextValue $wsQ4.Range("G3") "0.0133"
$wsQ4.Range("H3").Value = 1

# Row 4
Set-TextValue $wsQ4.Range("B4") "562530"
Set-TextValue $wsQ4.Range("C4") "华夏中证智选1000价值稳健策略ETF"
Set-TextValue $wsQ4.Range("D4") "0.36"
Set-TextValue $wsQ4.Range("E4") "96.22"
Set-TextValue $wsQ4.Range("F4") "0.95"
Set-TextValue $wsQ4.Range("G4") "0.0034"
$wsQ4.Range("H4").Value = 3

# Row 5
Set-TextValue $wsQ4.Range("B5") "519222"
Set-TextValue $wsQ4.Range("C5") "海富通欣益灵活配置混合A"
Set-TextValue $wsQ4.Range("D5") "0.25"
Set-TextValue $wsQ4.Range("E5") "31.65"
Set-TextValue $wsQ4.Range("F5") "0.20"
Set-TextValue $wsQ4.Range("G5") "0.0005"
$wsQ4.Range("H5").Value = 1

# Row 6
Set-TextValue $wsQ4.Range("B6") "519221"
Set-TextValue $wsQ4.Range("C6") "海富通欣益灵活配置混合C"
Set-TextValue $wsQ4.Range("D6") "0.10"
Set-TextValue $wsQ4.Range("E6") "31.65"
Set-TextValue $wsQ4.Range("F6") "0.20"
Set-TextValue $wsQ4.Range("G6") "0.0002"
$wsQ4.Range("H6").Value = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert the new 2022-Q4 summary row
#    at the top of the data (row 2) and push the rest down, adding a new
#    last row for 2020-Q4.
# ---------------------------------------------------------------------
$wsTotal.Range("A10").Value = 8
$wsTotal.Range("B10").Value = "2020-Q4"
$wsTotal.Range("C10").Value = 6
$wsTotal.Range("D10").Value = 0.76

$wsTotal.Range("A9").Value = 7
$wsTotal.Range("B9").Value = "2021-Q1"
$wsTotal.Range("C9").Value = 12
$wsTotal.Range("D9").Value = 0.97

$wsTotal.Range("A8").Value = 6
$wsTotal.Range("B8").Value = "2021-Q2"
$wsTotal.Range("C8").Value = 9
$wsTotal.Range("D8").Value = 0.7

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2021-Q3"
$wsTotal.Range("C7").Value = 5
$wsTotal.Range("D7").Value = 0.38

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q4"
$wsTotal.Range("C6").Value = 3
$wsTotal.Range("D6").Value = 4.34

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2022-Q1"
$wsTotal.Range("C5").Value = 3
$wsTotal.Range("D5").Value = 0.36

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.36

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 0.42

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.36

# ---------------------------------------------------------------------
# 3. Restore the originally-active/selected tab ("2020-Q4") — copying a
#    sheet makes the new copy active, so re-select the last tab to match
#    the source workbook's selection state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
